$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.172755002975464
$ws.Range("B1").Value = 2.43589186668396
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.366260290145874
$ws.Range("E1").Value = 1.23726224899292
